# SizingCalculator.xlsx - "fixed min vars in env; added more comments on velero items"
#
# The DJANGO_MEMORY_MAX / DJANGODB_MEMORY_MAX / VECTORDB_MEMORY_MAX "min"
# memory sizing inputs (column D, rows 12/15/16) were over-stated at 4Gi;
# correct them down to 3Gi on both the DEV and UAT & PROD sizing sheets.
# All of the dependent formula cells (H/I/J/K/L/M/N/O columns and the
# generated DJANGO_MEMORY_MAX=.. / DJANGODB_MEMORY_MAX=.. / VECTORDB_MEMORY_MAX=..
# env-var strings) recalculate automatically.

$wb = $excel.ActiveWorkbook

$devSheet = $wb.Worksheets.Item("DEV")
$uatSheet = $wb.Worksheets.Item("UAT & PROD")

# --- DEV sheet: minimum memory sizing fixes ---
$devSheet.Range("D12").Value = 3
$devSheet.Range("D15").Value = 3
$devSheet.Range("D16").Value = 3

# --- UAT & PROD sheet: same minimum memory sizing fixes ---
$uatSheet.Range("D12").Value = 3
$uatSheet.Range("D15").Value = 3
$uatSheet.Range("D16").Value = 3

# Leave the cursor parked on UAT & PROD!D16 (the last value we touched there)
# before switching back to DEV, which ends up the active/front sheet.
$uatSheet.Activate() | Out-Null
$uatSheet.Range("D16").Select() | Out-Null

$devSheet.Activate() | Out-Null
$devSheet.Range("D9").Select() | Out-Null
